$d = $word.ActiveDocument

$replacements = @(
    @("291×9=", "542×8="),
    @("247×8=", "785×5="),
    @("714×2=", "661×9="),
    @("977×3=", "665×6="),
    @("218×6=", "102×7="),
    @("612×5=", "782×4="),
    @("683×2=", "518×5="),
    @("560×4=", "353×5="),
    @("767×9=", "595×7="),
    @("838×6=", "588×9="),
    @("865×8=", "245×5="),
    @("113×8=", "436×8="),
    @("424×7=", "112×5="),
    @("670×5=", "866×9="),
    @("244×2=", "457×4="),
    @("715×8=", "481×7="),
    @("909×8=", "729×2="),
    @("571×3=", "589×9="),
    @("746×6=", "689×4="),
    @("546×9=", "772×6="),
    @("169×3=", "262×7="),
    @("153×3=", "312×9="),
    @("660×2=", "770×6="),
    @("576×4=", "190×5="),
    @("478×3=", "666×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
